$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Insert-ParagraphXmlAfter($range, [string]$innerBody) {
    # Inserts a brand-new paragraph (with the given inner <w:p>...</w:p> body)
    # immediately after the given (collapsed or not) range.
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="' + $wNs + '"><w:body>' + $innerBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

function Replace-ParagraphRuns($paragraph, [string]$innerRunsXml, [string]$plainTargetText) {
    # Rebuilds the *start* of a paragraph with explicit run XML (preserving the
    # paragraph's own <w:p ...> attributes), then deletes the old trailing
    # content that is now duplicated after the freshly-inserted text.
    $startPos = $paragraph.Range.Start
    $insertPoint = $d.Range($startPos, $startPos)
    $innerBody = '<w:p>' + $innerRunsXml + '</w:p>'
    Insert-ParagraphXmlAfter $insertPoint $innerBody

    $p2 = $d.Paragraphs($paragraph.Index)
    $newLen = $plainTargetText.Length
    $oldStart = $p2.Range.Start + $newLen
    $oldEnd = $p2.Range.End - 1
    if ($oldEnd -gt $oldStart) {
        $toDelete = $d.Range($oldStart, $oldEnd)
        $toDelete.Delete()
    }
}

# ---------------------------------------------------------------------------
# Change 1: new paragraph about text styles / color palette, right after the
# "(on rented)" paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(on rented)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterRented = $rng.End
$insertPoint = $d.Range($afterRented, $afterRented)
$newParaBody = '<w:p>' +
    '<w:r><w:t>When color palette is altered, text styles will be updated</w:t></w:r>' +
    '<w:r><w:t>. Maybe put in functionality that allows text styles to be bound to observables.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Hell maybe make them observables?</w:t></w:r>' +
    '</w:p>'
Insert-ParagraphXmlAfter $insertPoint $newParaBody

# ---------------------------------------------------------------------------
# Change 2: "Typed Pool: ElementPool<Element>" -- collapse the ": ",
# "ElementPool" and "<Element>" runs (and proofErr wrappers) into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Typed Pool", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$typedPoolPara = $rng.Paragraphs(1)
$runsXml = '<w:r w:rsidRPr="00905593"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Typed Pool</w:t></w:r>' +
           '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>: ElementPool&lt;Element&gt;</w:t></w:r>'
Replace-ParagraphRuns $typedPoolPara $runsXml "Typed Pool: ElementPool<Element>"

# ---------------------------------------------------------------------------
# Change 3: "The rect transform of the currently highlighted element..."
# -- collapse into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2 = $d.Content
$rng2.Find.Execute("rect transform of the currently highlighted", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rectPara = $rng2.Paragraphs(1)
$runsXml = '<w:r><w:t>The rect transform of the currently highlighted element (if there is one)</w:t></w:r>'
Replace-ParagraphRuns $rectPara $runsXml "The rect transform of the currently highlighted element (if there is one)"

# ---------------------------------------------------------------------------
# Change 4: "Boolean for flagging when the ui element is still moving"
# -- collapse into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Boolean for flagging", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boolPara = $rng.Paragraphs(1)
$runsXml = '<w:r><w:t>Boolean for flagging when the ui element is still moving</w:t></w:r>'
Replace-ParagraphRuns $boolPara $runsXml "Boolean for flagging when the ui element is still moving"

# ---------------------------------------------------------------------------
# Change 5: add <w:lastRenderedPageBreak/> to the start of the run containing
# "Drop shadow for Stats like Attack".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Drop shadow for Stats like Attack", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dropPara = $rng.Paragraphs(1)
$runsXml = '<w:r><w:lastRenderedPageBreak/><w:t>Drop shadow for Stats like Attack</w:t></w:r>' +
           '<w:r w:rsidR="00C95236"><w:t>, defence, health(?), magic</w:t></w:r>'
Replace-ParagraphRuns $dropPara $runsXml "Drop shadow for Stats like Attack, defence, health(?), magic"

Write-Output "All edits applied"
